# modelo_person.xlsx edit:
#  - fix header typo "bithday" -> "birthday"
#  - move the active cell selection from T9 to G7
#  - nudge a couple of column widths by a hundredth of a character

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in the header row (column F, "bithday" -> "birthday")
$ws.Cells.Replace("bithday", "birthday")

# Small column-width tweaks (column D and column Q)
$ws.Columns.Item(4).ColumnWidth = 3.79
$ws.Columns.Item(17).ColumnWidth = 7.85

# Update the current selection to G7
$ws.Range("G7").Select()
